$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.81714707784717
$ws.Range("L2").Value = 5.60241266790566
$ws.Range("B3").Value = 6.38308258262284
$ws.Range("L3").Value = 6.58293485463342
$ws.Range("B4").Value = 6.01644071117871
$ws.Range("B5").Value = 6.56604977851573
$ws.Range("L5").Value = 7.04790297162438
